$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new year (2022) column has been added to the right of the existing
# data (column P = 2021), so copy column P's formatting into the new
# column Q via an insert-shift (mirrors "copy column, insert copied
# cells" in the Excel UI) and then fill in the new values.
$ws.Range("P4:P5").Copy()
$ws.Range("Q4:Q5").Insert(-4161)

$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 64.2

# Move/record the active selection as it ended up after the edit.
$ws.Range("R4").Select()
